$wb = $excel.ActiveWorkbook

# Add the new "cost" sheet after the last existing sheet ("demand"),
# mirroring the author adding a 4th tab to the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "cost"

# Populate in the same left-to-right / top-to-bottom order the author
# used, so new shared strings land in the same sequence as the source
# workbook (hubs location variable, initial routing..., total operation
# cost of one hub, [0,1], $/Month, operation cost).
$ws.Range("A1").Value = "hubs location variable"
$ws.Range("A2").Value = "initial routing and consolidation efficiency"
$ws.Range("A3").Value = "total operation cost of one hub"
$ws.Range("D1").Value = "[0,1]"
$ws.Range("C3").Value = "$/Month"
$ws.Range("A4").Value = "operation cost"

$ws.Range("B1").Value = 0.7
$ws.Range("C1").Value = "Dmnl"

$ws.Range("B2").Value = 0.6
$ws.Range("C2").Value = "Dmnl"
$ws.Range("D2").Value = "[0,1]"

$ws.Range("B3").Value = 100

$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = "$/Month"

# Match the column A width used on the other sheets (best-fit style).
$ws.Columns.Item(1).ColumnWidth = 35.75

# Author left the selection on B9 on the new "cost" sheet, which also
# becomes the active tab (sheet index 4, activeTab goes 2 -> 3).
$ws.Range("B9").Select()
